# fix: Arvore de decisões corrigida
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TOWERS column (B) decision values corrected
$ws.Range("B2").Value = "NORMAL_T"
$ws.Range("B3").Value = "NORMAL_T"
$ws.Range("B4").Value = "BIG_T"

# MELIANTS column (C) decision values corrected
$ws.Range("C2").Value = "NORMAL_M"
$ws.Range("C3").Value = "NONE_M"
$ws.Range("C4").Value = "BIG_M"

# POSITION column (F) value corrected
$ws.Range("F2").Value = "5x4.20"

# Update selection to reflect the active cell in the saved view
$ws.Range("C4").Select()
